# ------------------------------------------------------------------
# Applies the "PO Forecast" update:
#  1. Rename header B1 on "Weekly Quantity"  -> "Weekly_PO_Qty"
#  2. Rename header B1 on "Monthly Trend"    -> "Monthly_PO_Qty"
#  3. Add a new worksheet "PO Forecast" (after "Monthly Trend") with
#     columns ds / PO_Forecast / yhat_lower / yhat_upper and the
#     forecast data rows.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header rename --------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header rename -----------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" worksheet at the end ----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$header = $wsForecast.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45340.99999999999, 39, 7.961621341697001, 71.75586530405224),
    @(45347.99999999999, 39, 10.19746556306054, 74.55260551194608),
    @(45361.99999999999, 37, 5.012134386100799, 69.54879722561734),
    @(45368.99999999999, 36, 3.964239245472245, 69.19853327395971),
    @(45375.99999999999, 36, 6.86942921611727, 67.18229861421236),
    @(45382.99999999999, 35, 3.277888892415623, 68.18653463292161),
    @(45389.99999999999, 34, 4.151753371281339, 66.02874060705973),
    @(45396.99999999999, 33, 3.919862250932204, 65.63645674195377),
    @(45403.99999999999, 33, 0.2933242043822515, 62.56879649014979),
    @(45410.99999999999, 32, -0.4774695139083556, 62.57844851296599),
    @(45424.99999999999, 30, 0.4551030050587972, 63.41413409850946),
    @(45431.99999999999, 30, -1.036734291383935, 58.88624131466359),
    @(45522.99999999999, 20, -13.71819345348438, 50.92854556601289),
    @(45564.99999999999, 15, -13.86788093347226, 46.6915502287987),
    @(45585.99999999999, 13, -16.35701204240086, 42.57472849629418),
    @(45592.99999999999, 12, -19.69258931204238, 45.22376805499871),
    @(45599.99999999999, 12, -19.92727024684887, 45.1952457053766),
    @(45606.99999999999, 11, -22.31229175397433, 42.45611464704582),
    @(45613.99999999999, 10, -22.89201665495384, 40.49027370102974),
    @(45620.99999999999, 9, -20.9264159636544, 41.39084655200069),
    @(45627.99999999999, 9, -22.01121751955069, 39.93959783216765),
    @(45634.99999999999, 8, -24.12290462415216, 41.41813845980244),
    @(45641.99999999999, 7, -26.11847019435924, 39.53592937130932)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "PO Forecast sheet created with" $data.Length "rows"
